$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.74
$ws.Range("A10").Value = -21.552
$ws.Range("A12").Value = -21.58
$ws.Range("C12").Value = -10.563
$ws.Range("C17").Value = -13.36
$ws.Range("A18").Value = -21.895
$ws.Range("C26").Value = -13.313
$ws.Range("C27").Value = -13.657
$ws.Range("C28").Value = -13.076
$ws.Range("A37").Value = -19.909
$ws.Range("C37").Value = -12.813
$ws.Range("A55").Value = -21.795
$ws.Range("C65").Value = -12.1
$ws.Range("A68").Value = -21.536
$ws.Range("C73").Value = -12.264
$ws.Range("A77").Value = -20.651
$ws.Range("A78").Value = -20.107
$ws.Range("C84").Value = -13.947
$ws.Range("C85").Value = -12.082
$ws.Range("C93").Value = -10.901
$ws.Range("C95").Value = -11.682
$ws.Range("C98").Value = -13.105
$ws.Range("C99").Value = -11.312
$ws.Range("C101").Value = -12.37
